$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions automated refresh)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.892.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.917.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.48"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.84"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000225"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.63"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.405.35"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.971.37"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.69"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.922.15"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "429.48"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.38"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.681"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.07"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.40"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.92"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.88%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.42"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.107"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0846"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.62"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.02"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.124"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.98"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.57"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.286"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.33"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "373.80"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0346"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.718.23"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.62"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.94"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.96%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.125"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.79%  "

Write-Host "Updated cryptos list"